$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price/volume figures (and the PEPE/Monero/Stellar row reorder) for the
# Tue Apr 16 14:42:25 UTC 2024 cryptos-list refresh. Each entry below is the new
# B/C/D/E content for that sheet row; only cells that actually changed are listed.
$updates = @(
    @{ Row=2; D="62.271.47"; E="  -4.34%  " }
    @{ Row=3; D="3.050.29"; E="  -3.98%  " }
    @{ Row=4; E="  -0.20%  " }
    @{ Row=5; D="533.73"; E="  -5.90%  " }
    @{ Row=6; D="131.39"; E="  -10.89%  " }
    @{ Row=7; E="  +0.07%  " }
    @{ Row=8; D="3.039.80"; E="  -4.07%  " }
    @{ Row=9; D="0.484"; E="  -3.74%  " }
    @{ Row=10; D="0.151"; E="  -4.51%  " }
    @{ Row=11; D="6.06"; E="  -11.49%  " }
    @{ Row=12; D="0.450"; E="  -5.42%  " }
    @{ Row=13; D="0.0000221"; E="  -2.48%  " }
    @{ Row=14; D="33.89"; E="  -9.96%  " }
    @{ Row=15; D="3.510.59"; E="  -5.22%  " }
    @{ Row=16; D="62.293.30"; E="  -4.64%  " }
    @{ Row=17; E="  -3.10%  " }
    @{ Row=18; D="3.067.23"; E="  -3.91%  " }
    @{ Row=19; D="6.49"; E="  -6.87%  " }
    @{ Row=20; D="472.06"; E="  -10.30%  " }
    @{ Row=21; D="13.12"; E="  -8.13%  " }
    @{ Row=22; D="0.688"; E="  -5.55%  " }
    @{ Row=23; D="7.06"; E="  -7.14%  " }
    @{ Row=24; D="77.90"; E="  -2.40%  " }
    @{ Row=25; D="11.87"; E="  -9.47%  " }
    @{ Row=26; D="0.997"; E="  +0.06%  " }
    @{ Row=27; D="2.66"; E="  -7.71%  " }
    @{ Row=28; D="8.12"; E="  -11.19%  " }
    @{ Row=29; E="  -0.05%  " }
    @{ Row=30; D="25.43"; E="  -5.68%  " }
    @{ Row=31; D="1.83"; E="  -16.95%  " }
    @{ Row=32; E="  -5.89%  " }
    @{ Row=33; D="2.36"; E="  -10.95%  " }
    @{ Row=34; D="56.62"; E="  +3.82%  " }
    @{ Row=35; D="5.84"; E="  -5.89%  " }
    @{ Row=36; D="5.13"; E="  -6.52%  " }
    @{ Row=37; D="468.31"; E="  -15.02%  " }
    @{ Row=38; D="3.093.42"; E="  -2.75%  " }
    @{ Row=39; D="0.0386"; E="  -12.27%  " }
    @{ Row=40; D="0.0780"; E="  -7.07%  " }
    @{ Row=41; D="7.93"; E="  -6.02%  " }
    @{ Row=42; D="0.111"; E="  -11.88%  " }
    @{ Row=43; D="2.55"; E="  -9.68%  " }
    @{ Row=44; E="  +0.04%  " }
    @{ Row=45; D="0.245"; E="  -10.43%  " }
    @{ Row=46; D="1.99"; E="  -11.77%  " }
    @{ Row=47; D="23.98" }
    @{ Row=48; B="PEPE"; C="https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D="0.0₃0511"; E="  -4.71%  " }
    @{ Row=49; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="116.43"; E="  -5.66%  " }
    @{ Row=50; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.106"; E="  -4.03%  " }
    @{ Row=51; D="1.96"; E="  -8.98%  " }
)

foreach ($u in $updates) {
    $row = $u.Row
    foreach ($col in "B", "C", "D", "E") {
        if (-not $u.ContainsKey($col)) { continue }
        $cell = $ws.Range("$col$row")
        if ($col -eq "D" -or $col -eq "E") {
            # Price/volume columns are stored as text (e.g. "62.271.47", "  -4.34%  ");
            # force text format so Excel does not reinterpret them as numbers, then
            # restore the default style so no stray formatting is left behind.
            $cell.NumberFormat = "@"
            $cell.Value = $u[$col]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $u[$col]
        }
    }
}
